$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('1amag5h', 'Tried acrylics for the first time!', 'https://i.redd.it/c94znyipbghc1.jpeg')
    ,@('1am4fp0', 'What step am I missing with my gel?', 'https://www.reddit.com/r/Nails/comments/1am4fp0/what_step_am_i_missing_with_my_gel/')
    ,@('1alw12k', 'Moons and stars 🌚', 'https://i.redd.it/4ngiiuu88dhc1.jpeg')
    ,@('1alqhtk', 'Practising my hand painting', 'https://i.redd.it/d4ge8omrgbhc1.jpeg')
    ,@('1amgdyw', 'LynBDesigns Sale?', 'https://www.reddit.com/r/RedditLaqueristas/comments/1amgdyw/lynbdesigns_sale/')
    ,@('1amg9ka', 'More LynB Designs ~ One on each hand!', 'https://www.reddit.com/gallery/1amg9ka')
    ,@('1amg96t', 'Steel Your Heart (first attempt mixing my own polish)', 'https://www.reddit.com/gallery/1amg96t')
    ,@('1amedws', 'INLP Eclipse is just a great color', 'https://www.reddit.com/gallery/1amedws')
    ,@('1amcx63', 'Felt compelled to paint my nails red to celebrate Lunar New Year / the Suberb Owl / Valentine’s Day ❣️', 'https://www.reddit.com/gallery/1amcx63')
    ,@('1amcgf0', 'Tết Nails!', 'https://www.reddit.com/gallery/1amcgf0')
    ,@('1ambxx4', 'Cadillacquer Winter Collection Pre Order?', 'https://www.reddit.com/r/RedditLaqueristas/comments/1ambxx4/cadillacquer_winter_collection_pre_order/')
    ,@('1amanva', 'I love it when my nails match things 🩵', 'https://www.reddit.com/gallery/1amanva')
    ,@('1ama7cc', 'ISO help to not rip off my nails while sleeping 🥲', 'https://www.reddit.com/gallery/1ama7cc')
    ,@('1am9zw5', 'Kind of cool/weird ILNP combo', 'https://www.reddit.com/gallery/1am9zw5')
    ,@('1am9jjw', 'Calling this one “Blood and Chocolate” ❤️🥀🍫🤎', 'https://www.reddit.com/gallery/1am9jjw')
    ,@('1am8810', 'KBS - Lights Out', 'https://v.redd.it/7t0uzqyytfhc1')
    ,@('1am87c3', 'we kinda match', 'https://i.redd.it/3pnhrq7vtfhc1.jpeg')
    ,@('1am5q5a', 'Thermal Valentine’s Nails', 'https://i.redd.it/wqnih4l6bfhc1.jpeg')
    ,@('1am7rpy', 'Tenta cool vs Don''t Provoke Me', 'https://www.reddit.com/r/RedditLaqueristas/comments/1am7rpy/tenta_cool_vs_dont_provoke_me/')
    ,@('1am7p8b', 'Nothing says valentine''s month like royal blue China 💙🍽️', 'https://www.reddit.com/gallery/1am7p8b')
    ,@('1am7ir1', 'A student told me my nails look like “an alien threw up on them” 👽', 'https://www.reddit.com/gallery/1am7ir1')
    ,@('1am7h56', 'New to nail obsession, some questions', 'https://www.reddit.com/r/RedditLaqueristas/comments/1am7h56/new_to_nail_obsession_some_questions/')
    ,@('1am76ku', 'Kind of vday inspired', 'https://i.redd.it/j2ss6rk5mfhc1.jpeg')
    ,@('1am6ojc', 'First try with super cheap press-ons', 'https://i.redd.it/oub2wghcifhc1.jpeg')
    ,@('1am6ntn', 'Rosewater Jelly with a lil something Sassy', 'https://www.reddit.com/gallery/1am6ntn')
    ,@('1am6i7o', 'Tried nail art for the first time… it was much cuter in my head', 'https://www.reddit.com/gallery/1am6i7o')
    ,@('1am68lu', 'Spellbound Nails', 'https://www.reddit.com/r/RedditLaqueristas/comments/1am68lu/spellbound_nails/')
    ,@('1am658t', 'Valentine''s Day skittle 🩷💜', 'https://i.redd.it/w7dy64xbefhc1.jpeg')
    ,@('1am5oib', 'Totally fan- girling....', 'https://i.redd.it/51cu2lmtafhc1.jpeg')
    ,@('1am54g7', 'Did anyone get the Zoya Love polish?', 'https://www.reddit.com/r/RedditLaqueristas/comments/1am54g7/did_anyone_get_the_zoya_love_polish/')
    ,@('1am4t3d', 'Any more opaque/matte "Boring" flakie toppers out there?', 'https://www.reddit.com/r/RedditLaqueristas/comments/1am4t3d/any_more_opaquematte_boring_flakie_toppers_out/')
    ,@('1am4atu', '🩷💙Neon Hearts💜💛', 'https://www.reddit.com/gallery/1am4atu')
    ,@('1am46vv', 'OPI - Kyoto Pearl', 'https://i.redd.it/fbg66sfqzehc1.jpeg')
    ,@('1am3t4a', 'Any reviews on Essie special effects line?', 'https://www.reddit.com/r/RedditLaqueristas/comments/1am3t4a/any_reviews_on_essie_special_effects_line/')
    ,@('1am3nuf', 'Subtle Valentines nails ❤️', 'https://www.reddit.com/gallery/1am3nuf')
    ,@('1am3b2f', 'Arctic Fox "Cosmic Sunshine" + Holo Taco "Lunar Unicorn Skin"', 'https://i.redd.it/h9w8z011tehc1.jpeg')
    ,@('1am2kzj', 'Multichrome Gradient', 'https://v.redd.it/seincgjnnehc1')
    ,@('1am2kmj', 'Lovestruck Collection Valentine’s Day mix n match mani 🥰 (PR - gifted)', 'https://www.reddit.com/gallery/1am2kmj')
    ,@('1am2gjn', 'silver nail polish has my heart', 'https://i.redd.it/5ksafodqmehc1.jpeg')
    ,@('1am2buq', 'ILNP “Love Language”', 'https://i.redd.it/2n3zfrerlehc1.jpeg')
    ,@('1am27t7', 'Palate cleanser :P', 'https://i.redd.it/xmvosw3ykehc1.jpeg')
    ,@('1am25m0', 'Any Mooncat Snake Charmer dupes out there?', 'https://www.reddit.com/r/RedditLaqueristas/comments/1am25m0/any_mooncat_snake_charmer_dupes_out_there/')
    ,@('1am21qu', 'First attempt at nail "art"', 'https://i.redd.it/orivcn8pjehc1.jpeg')
    ,@('1am1u3r', 'It’s beautiful but AT WHAT COST (also I would love recs for more vampy shades with red shimmer)', 'https://www.reddit.com/gallery/1am1u3r')
    ,@('1am0k4y', 'Wildflower Lacquer Dark Beauty (PPU Jan 2024)', 'https://www.reddit.com/gallery/1am0k4y')
    ,@('1alzott', 'Mooncat Shattered Glass', 'https://www.reddit.com/gallery/1alzott')
    ,@('1alyfyz', 'Barb from Sassy Sauce', 'https://i.redd.it/l4csvcsfsdhc1.jpeg')
    ,@('1alxgva', 'blue/brown duochrome polish?', 'https://i.redd.it/ge1yhchmkdhc1.jpeg')
    ,@('1alxeiv', 'Spirited Away Studio Ghibli Nails 👻🐉🐖🐀⚫', 'https://www.reddit.com/gallery/1alxeiv')
    ,@('1alxc6h', 'Polished for Days I Had Wings Once', 'https://i.redd.it/bfd4vb0kjdhc1.jpeg')
    ,@('1alx8z9', 'first time attempting french nails', 'https://www.reddit.com/gallery/1alx8z9')
    ,@('1alx6w8', 'The color is subtle, but that holo...', 'https://i.redd.it/zxh1mrnbidhc1.png')
    ,@('1alwvku', 'Sheer pink mani', 'https://i.redd.it/18doe83qfdhc1.png')
    ,@('1alp72m', 'Does anyone know how to clean up nail polish thinner?', 'https://www.reddit.com/r/RedditLaqueristas/comments/1alp72m/does_anyone_know_how_to_clean_up_nail_polish/')
    ,@('1alt2us', 'Can you tell pink is my favorite color? 🩷🌸🌷🎀 Polishes listed below', 'https://www.reddit.com/gallery/1alt2us')
    ,@('1alv0kb', 'Wall Reorganization Day!', 'https://i.redd.it/nbbo5fb3zchc1.jpeg')
    ,@('1amiahv', 'Valentines nails 🩷', 'https://www.reddit.com/gallery/1amiahv')
    ,@('1amg6jb', 'Silver heart nails 🩷🤍', 'https://www.reddit.com/gallery/1amg6jb')
    ,@('1amejtq', 'Hand painted', 'https://i.redd.it/3znzgwcebhhc1.jpeg')
    ,@('1amcbu4', 'Valentine’s Day nails 💅🏼', 'https://www.reddit.com/gallery/1amcbu4')
    ,@('1ambx1u', 'Fnaf nail art', 'https://www.reddit.com/gallery/1ambx1u')
    ,@('1amahji', 'i suck at making a decision🥲😅', 'https://www.reddit.com/gallery/1amahji')
    ,@('1am9fcs', 'Scooby Doo', 'https://i.redd.it/48798k5e3ghc1.jpeg')
    ,@('1am8mzf', 'did my valentine‘s day nails today🤍', 'https://v.redd.it/72di1oc7xfhc1')
    ,@('1am54z0', 'another vday set for your feed 🙈', 'https://www.reddit.com/gallery/1am54z0')
    ,@('1am4jmg', 'My hand painted press-on nails.', 'https://www.reddit.com/gallery/1am4jmg')
    ,@('1am141k', 'Cold but comfy', 'https://www.reddit.com/gallery/1am141k')
    ,@('1am109z', 'New set I did. Wearing a gold and green dress to a wedding this weekend and didn’t want traditional valentines colors to clash', 'https://www.reddit.com/gallery/1am109z')
    ,@('1am0bnh', 'Self taught', 'https://www.reddit.com/gallery/1am0bnh')
    ,@('1alyczq', 'Why is my nail foil peeling and scratching off? Issue with the glue?', 'https://www.reddit.com/r/NailArt/comments/1alyczq/why_is_my_nail_foil_peeling_and_scratching_off/')
    ,@('1alwyqj', 'Animal print for my client from earlier today', 'https://i.redd.it/deajo7xggdhc1.jpeg')
    ,@('1alwn9k', 'Moons and stars ✨', 'https://i.redd.it/p0x43yrrddhc1.jpeg')
)

$startRow = 713
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
